# --- Auto-generated COM-interop edit script -------------------------
# Adds the '2022-Q3' sheet (with fund-holding detail rows) right after
# '总计' (Total) and inserts the matching summary row in '总计'.

$wb = $excel.ActiveWorkbook

# 1) Insert the new worksheet before the current 2nd sheet (2022-Q1),
#    i.e. right after '总计', and name it 2022-Q3.
$insertBefore = $wb.Worksheets.Item(2)
$newSheet = $wb.Worksheets.Add($insertBefore)
$newSheet.Name = "2022-Q3"

# Sheet that already has the correct header/column-A formatting we want
# to mirror (the old 2022-Q1 sheet, now shifted one slot to the right).
$fmtSource = $wb.Worksheets.Item(3)

# Copy the bold header-row formatting (B1:H1).
$fmtSource.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

# Copy the column-A index-cell formatting down for all 21 data rows.
$fmtSource.Range("A2").Copy()
$newSheet.Range("A2:A22").PasteSpecial(-4122)

# Header labels
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Columns B, D, E, F, G hold numeric-looking text (fund codes with leading
# zeros, and decimal figures stored as strings in the source data) -- force
# them to Text before writing so Excel doesn't auto-coerce to numbers, then
# reset the number format back to Normal so no stray style survives.
$newSheet.Range("B2:B22").NumberFormat = "@"
$newSheet.Range("D2:D22").NumberFormat = "@"
$newSheet.Range("E2:E22").NumberFormat = "@"
$newSheet.Range("F2:F22").NumberFormat = "@"
$newSheet.Range("G2:G20").NumberFormat = "@"

# Row 2: 011346
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "011346"
$newSheet.Range("C2").Value = "淳厚鑫淳一年持有期混合"
$newSheet.Range("D2").Value = "4.81"
$newSheet.Range("E2").Value = "69.72"
$newSheet.Range("F2").Value = "4.11"
$newSheet.Range("G2").Value = "0.1977"
$newSheet.Range("H2").Value = 2
# Row 3: 000006
$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "000006"
$newSheet.Range("C3").Value = "西部利得量化成长混合A"
$newSheet.Range("D3").Value = "13.82"
$newSheet.Range("E3").Value = "86.11"
$newSheet.Range("F3").Value = "1.37"
$newSheet.Range("G3").Value = "0.1893"
$newSheet.Range("H3").Value = 4
# Row 4: 020015
$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").Value = "020015"
$newSheet.Range("C4").Value = "国泰区位优势混合A"
$newSheet.Range("D4").Value = "1.95"
$newSheet.Range("E4").Value = "87.12"
$newSheet.Range("F4").Value = "5.46"
$newSheet.Range("G4").Value = "0.1065"
$newSheet.Range("H4").Value = 4
# Row 5: 009874
$newSheet.Range("A5").Value = 3
$newSheet.Range("B5").Value = "009874"
$newSheet.Range("C5").Value = "九泰久睿量化股票A"
$newSheet.Range("D5").Value = "3.15"
$newSheet.Range("E5").Value = "93.83"
$newSheet.Range("F5").Value = "3.30"
$newSheet.Range("G5").Value = "0.1040"
$newSheet.Range("H5").Value = 7
# Row 6: 010779
$newSheet.Range("A6").Value = 4
$newSheet.Range("B6").Value = "010779"
$newSheet.Range("C6").Value = "西部利得量化优选一年持有期混合A"
$newSheet.Range("D6").Value = "4.88"
$newSheet.Range("E6").Value = "88.29"
$newSheet.Range("F6").Value = "2.08"
$newSheet.Range("G6").Value = "0.1015"
$newSheet.Range("H6").Value = 3
# Row 7: 012454
$newSheet.Range("A7").Value = 5
$newSheet.Range("B7").Value = "012454"
$newSheet.Range("C7").Value = "淳厚鑫悦混合A"
$newSheet.Range("D7").Value = "2.06"
$newSheet.Range("E7").Value = "75.61"
$newSheet.Range("F7").Value = "4.80"
$newSheet.Range("G7").Value = "0.0989"
$newSheet.Range("H7").Value = 2
# Row 8: 007126
$newSheet.Range("A8").Value = 6
$newSheet.Range("B8").Value = "007126"
$newSheet.Range("C8").Value = "博道远航混合A"
$newSheet.Range("D8").Value = "7.01"
$newSheet.Range("E8").Value = "88.04"
$newSheet.Range("F8").Value = "0.91"
$newSheet.Range("G8").Value = "0.0638"
$newSheet.Range("H8").Value = 9
# Row 9: 007127
$newSheet.Range("A9").Value = 7
$newSheet.Range("B9").Value = "007127"
$newSheet.Range("C9").Value = "博道远航混合C"
$newSheet.Range("D9").Value = "5.38"
$newSheet.Range("E9").Value = "88.04"
$newSheet.Range("F9").Value = "0.91"
$newSheet.Range("G9").Value = "0.0490"
$newSheet.Range("H9").Value = 9
# Row 10: 012455
$newSheet.Range("A10").Value = 8
$newSheet.Range("B10").Value = "012455"
$newSheet.Range("C10").Value = "淳厚鑫悦混合C"
$newSheet.Range("D10").Value = "0.68"
$newSheet.Range("E10").Value = "75.61"
$newSheet.Range("F10").Value = "4.80"
$newSheet.Range("G10").Value = "0.0326"
$newSheet.Range("H10").Value = 2
# Row 11: 010780
$newSheet.Range("A11").Value = 9
$newSheet.Range("B11").Value = "010780"
$newSheet.Range("C11").Value = "西部利得量化优选一年持有期混合C"
$newSheet.Range("D11").Value = "1.22"
$newSheet.Range("E11").Value = "88.29"
$newSheet.Range("F11").Value = "2.08"
$newSheet.Range("G11").Value = "0.0254"
$newSheet.Range("H11").Value = 3
# Row 12: 011228
$newSheet.Range("A12").Value = 10
$newSheet.Range("B12").Value = "011228"
$newSheet.Range("C12").Value = "西部利得量化成长混合C"
$newSheet.Range("D12").Value = "1.67"
$newSheet.Range("E12").Value = "86.11"
$newSheet.Range("F12").Value = "1.37"
$newSheet.Range("G12").Value = "0.0229"
$newSheet.Range("H12").Value = 4
# Row 13: 002137
$newSheet.Range("A13").Value = 11
$newSheet.Range("B13").Value = "002137"
$newSheet.Range("C13").Value = "诺安利鑫灵活配置混合A"
$newSheet.Range("D13").Value = "0.44"
$newSheet.Range("E13").Value = "76.46"
$newSheet.Range("F13").Value = "4.73"
$newSheet.Range("G13").Value = "0.0208"
$newSheet.Range("H13").Value = 1
# Row 14: 010120
$newSheet.Range("A14").Value = 12
$newSheet.Range("B14").Value = "010120"
$newSheet.Range("C14").Value = "九泰久福量化股票A"
$newSheet.Range("D14").Value = "0.54"
$newSheet.Range("E14").Value = "93.91"
$newSheet.Range("F14").Value = "3.31"
$newSheet.Range("G14").Value = "0.0179"
$newSheet.Range("H14").Value = 7
# Row 15: 001897
$newSheet.Range("A15").Value = 13
$newSheet.Range("B15").Value = "001897"
$newSheet.Range("C15").Value = "九泰久盛量化先锋灵活配置混合A"
$newSheet.Range("D15").Value = "0.50"
$newSheet.Range("E15").Value = "93.59"
$newSheet.Range("F15").Value = "3.29"
$newSheet.Range("G15").Value = "0.0164"
$newSheet.Range("H15").Value = 7
# Row 16: 009043
$newSheet.Range("A16").Value = 14
$newSheet.Range("B16").Value = "009043"
$newSheet.Range("C16").Value = "九泰久信量化股票"
$newSheet.Range("D16").Value = "0.43"
$newSheet.Range("E16").Value = "93.60"
$newSheet.Range("F16").Value = "3.31"
$newSheet.Range("G16").Value = "0.0142"
$newSheet.Range("H16").Value = 7
# Row 17: 004510
$newSheet.Range("A17").Value = 15
$newSheet.Range("B17").Value = "004510"
$newSheet.Range("C17").Value = "九泰久盛量化先锋灵活配置混合C"
$newSheet.Range("D17").Value = "0.28"
$newSheet.Range("E17").Value = "93.59"
$newSheet.Range("F17").Value = "3.29"
$newSheet.Range("G17").Value = "0.0092"
$newSheet.Range("H17").Value = 7
# Row 18: 014521
$newSheet.Range("A18").Value = 16
$newSheet.Range("B18").Value = "014521"
$newSheet.Range("C18").Value = "诺安利鑫灵活配置混合C"
$newSheet.Range("D18").Value = "0.05"
$newSheet.Range("E18").Value = "76.46"
$newSheet.Range("F18").Value = "4.73"
$newSheet.Range("G18").Value = "0.0024"
$newSheet.Range("H18").Value = 1
# Row 19: 007808
$newSheet.Range("A19").Value = 17
$newSheet.Range("B19").Value = "007808"
$newSheet.Range("C19").Value = "北信瑞丰量化优选灵活配置混合"
$newSheet.Range("D19").Value = "0.17"
$newSheet.Range("E19").Value = "78.75"
$newSheet.Range("F19").Value = "1.12"
$newSheet.Range("G19").Value = "0.0019"
$newSheet.Range("H19").Value = 5
# Row 20: 010121
$newSheet.Range("A20").Value = 18
$newSheet.Range("B20").Value = "010121"
$newSheet.Range("C20").Value = "九泰久福量化股票C"
$newSheet.Range("D20").Value = "0.04"
$newSheet.Range("E20").Value = "93.91"
$newSheet.Range("F20").Value = "3.31"
$newSheet.Range("G20").Value = "0.0013"
$newSheet.Range("H20").Value = 7
# Row 21: 015594
$newSheet.Range("A21").Value = 19
$newSheet.Range("B21").Value = "015594"
$newSheet.Range("C21").Value = "国泰区位优势混合C"
$newSheet.Range("D21").Value = "0.00"
$newSheet.Range("E21").Value = "87.12"
$newSheet.Range("F21").Value = "5.46"
$newSheet.Range("G21").Value = 0
$newSheet.Range("H21").Value = 4
# Row 22: 016399
$newSheet.Range("A22").Value = 20
$newSheet.Range("B22").Value = "016399"
$newSheet.Range("C22").Value = "九泰久睿量化股票C"
$newSheet.Range("D22").Value = "0.00"
$newSheet.Range("E22").Value = "93.83"
$newSheet.Range("F22").Value = "3.30"
$newSheet.Range("G22").Value = 0
$newSheet.Range("H22").Value = 7

# Drop back to the Normal style so the forced-text cells don't keep the
# temporary "@" number format in the saved file.
$newSheet.Range("B2:B22").Style = "Normal"
$newSheet.Range("D2:D22").Style = "Normal"
$newSheet.Range("E2:E22").Style = "Normal"
$newSheet.Range("F2:F22").Style = "Normal"
$newSheet.Range("G2:G20").Style = "Normal"

# 2) Update the '总计' (Total) summary sheet: insert a new row 2 for
#    2022-Q3 and push the existing rows down (values are unchanged, only
#    their row position moves).
$total = $wb.Worksheets.Item(1)
$total.Rows.Item(2).Insert()

# Re-apply the column-A index style + row formatting that Insert() does
# not fully carry over, by copying it from the row just below (row 3,
# which still has the original formatting).
$total.Range("A3:D3").Copy()
$total.Range("A2:D2").PasteSpecial(-4122)

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 21
$total.Range("D2").Value = 1.08
